$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 14 (budoyoni2@gmail.com review). Rows below shift up,
# so the former row 15 (ronoren61@gmail.com / nitanoren23@gmail.com / "amazing
# series of app...") becomes the new row 14.
$ws.Rows.Item(14).Delete()

$ws.Range("A14").Select()
